# Refresh the cryptocurrency price/volume table (cryptos.xlsx)
# "Updated cryptos list on Wed Oct 11 22:08:53 UTC 2023 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.754.55'
$ws.Range('E2').Value = '  -2.60%  '
$ws.Range('D3').Value = '1.567.26'
$ws.Range('E3').Value = '  -0.04%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '206.46'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').Value = '''0.490'
$ws.Range('E6').Value = '  -2.37%  '
$ws.Range('D8').Value = '21.91'
$ws.Range('E8').Value = '  -1.18%  '
$ws.Range('E9').Value = '  -0.91%  '
$ws.Range('E10').Value = '  -1.41%  '
$ws.Range('E11').Value = '  -0.32%  '
$ws.Range('D12').Value = '1.789.81'
$ws.Range('E12').Value = '  +0.01%  '
$ws.Range('D13').Value = '1.560.65'
$ws.Range('E13').Value = '  -0.95%  '
$ws.Range('D14').Value = '3.73'
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').Value = '0.514'
$ws.Range('E15').Value = '  -0.98%  '
$ws.Range('D16').Value = '26.788.32'
$ws.Range('E16').Value = '  -2.43%  '
$ws.Range('D17').Value = '61.35'
$ws.Range('E17').Value = '  -3.62%  '
$ws.Range('B18').Value = 'Chainlink'
$ws.Range('C18').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D18').Value = '''7.40'
$ws.Range('E18').Value = '  +1.86%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '215.06'
$ws.Range('E19').Value = '  +0.47%  '
$ws.Range('D20').Value = '0.0₃0677'
$ws.Range('E20').Value = '  -1.94%  '
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').Value = '''4.10'
$ws.Range('E22').Value = '  -0.21%  '
$ws.Range('D23').Value = '''9.30'
$ws.Range('E23').Value = '  -2.73%  '
$ws.Range('E24').Value = '  -0.93%  '
$ws.Range('E25').Value = '  -0.33%  '
$ws.Range('D26').Value = '6.73'
$ws.Range('E26').Value = '  +0.44%  '
$ws.Range('D27').Value = '14.95'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  -1.56%  '
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('E31').Value = '  -3.47%  '
$ws.Range('D32').Value = '3.16'
$ws.Range('E32').Value = '  -1.44%  '
$ws.Range('D33').Value = '1.394.73'
$ws.Range('E33').Value = '  +1.29%  '
$ws.Range('E34').Value = '  -1.38%  '
$ws.Range('E35').Value = '  -0.53%  '
$ws.Range('E36').Value = '  -1.08%  '
$ws.Range('D37').Value = '0.931'
$ws.Range('E37').Value = '  -2.88%  '
$ws.Range('D38').Value = '0.0163'
$ws.Range('E38').Value = '  -2.82%  '
$ws.Range('D39').Value = '0.529'
$ws.Range('E39').Value = '  -4.07%  '
$ws.Range('D40').Value = '0.819'
$ws.Range('E40').Value = '  -0.88%  '
$ws.Range('E41').Value = '  +0.06%  '
$ws.Range('D42').Value = '0.987'
$ws.Range('E42').Value = '  +0.89%  '
$ws.Range('D43').Value = '''1.80'
$ws.Range('E43').Value = '  -0.24%  '
$ws.Range('D44').Value = '5.34'
$ws.Range('E44').Value = '  +1.52%  '
$ws.Range('E45').Value = '  +0.76%  '
$ws.Range('D46').Value = '63.22'
$ws.Range('E46').Value = '  -1.73%  '
$ws.Range('D47').Value = '1.702.42'
$ws.Range('E47').Value = '  +0.14%  '
$ws.Range('D48').Value = '85.91'
$ws.Range('E48').Value = '  +0.67%  '
$ws.Range('D49').Value = '0.0₇0984'
$ws.Range('E49').Value = '  -1.88%  '
$ws.Range('D50').Value = '0.0952'
$ws.Range('E50').Value = '  -0.65%  '
$ws.Range('D51').Value = '0.0492'
$ws.Range('E51').Value = '  -0.83%  '
